$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.331.15"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.877.75"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7234"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.89"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08019"
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3157"
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.02"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08214"
$ws.Range("E11").Value = "  -2.19%  "
$ws.Range("D12").Value = "1.883.00"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.73"
$ws.Range("E13").Value = "  +4.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.224"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7126"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.416"
$ws.Range("E16").Value = "  +5.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008516"
$ws.Range("E17").Value = "  +3.98%  "
$ws.Range("D18").Value = "29.334.88"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.43"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.28"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.773"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1605"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.58"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.042"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.52"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.409"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.185"
$ws.Range("E31").Value = "  -8.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05366"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7610"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.177"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.705"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01876"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "1.274.98"
$ws.Range("E38").Value = "  +3.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.752"
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.436"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "113.05"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9094"
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "74.38"
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("E44").Value = "  +8.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "2.030.37"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5227"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.794"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.496"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4349"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.107"
$ws.Range("E51").Value = "  +0.59%  "
